# Apply the cryptos-list price/volume refresh described by the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row updates: column D (Price, text) and column E (Volume 1h, text).
# $null means that column is unchanged for that row.
$updates = @(
    @{ Row = 2; D = "63.308.53"; E = "  +0.78%  " }
    @{ Row = 3; D = "3.089.95"; E = "  -0.56%  " }
    @{ Row = 5; D = "582.48"; E = "  -0.48%  " }
    @{ Row = 6; D = "144.75"; E = "  +0.53%  " }
    @{ Row = 8; D = "3.083.29"; E = "  -0.56%  " }
    @{ Row = 9; D = "0.526"; E = "  -0.48%  " }
    @{ Row = 10; D = "0.158"; E = "  +6.16%  " }
    @{ Row = 11; D = "5.62"; E = "  -1.56%  " }
    @{ Row = 12; D = "0.454"; E = "  -2.82%  " }
    @{ Row = 13; D = "0.0000245"; E = "  +0.21%  " }
    @{ Row = 14; D = "37.36"; E = "  +5.55%  " }
    @{ Row = 15; D = $null; E = "  -1.09%  " }
    @{ Row = 16; D = "3.612.01"; E = "  -0.28%  " }
    @{ Row = 17; D = "63.274.74"; E = "  +0.84%  " }
    @{ Row = 18; D = "7.10"; E = "  -1.35%  " }
    @{ Row = 19; D = "3.094.85"; E = "  -0.31%  " }
    @{ Row = 20; D = "458.50"; E = "  -0.92%  " }
    @{ Row = 21; D = "14.18"; E = "  +0.88%  " }
    @{ Row = 22; D = "0.723"; E = "  -0.73%  " }
    @{ Row = 23; D = "7.42"; E = "  -1.55%  " }
    @{ Row = 24; D = "12.94"; E = "  -3.08%  " }
    @{ Row = 25; D = "80.95"; E = "  -1.57%  " }
    @{ Row = 26; D = "2.11"; E = "  -2.43%  " }
    @{ Row = 27; D = $null; E = "  +0.00%  " }
    @{ Row = 28; D = "8.89"; E = "  +7.52%  " }
    @{ Row = 29; D = $null; E = "  +0.04%  " }
    @{ Row = 30; D = "2.66"; E = "  -0.37%  " }
    @{ Row = 31; D = $null; E = "  -1.74%  " }
    @{ Row = 32; D = "6.77"; E = "  -0.78%  " }
    @{ Row = 33; D = "26.64"; E = "  -0.97%  " }
    @{ Row = 34; D = $null; E = "  -2.14%  " }
    @{ Row = 35; D = "0.0₃0841"; E = "  +2.11%  " }
    @{ Row = 36; D = $null; E = "  -1.21%  " }
    @{ Row = 37; D = "2.30"; E = "  -2.90%  " }
    @{ Row = 38; D = "3.33"; E = "  +6.37%  " }
    @{ Row = 39; D = "6.00"; E = "  -0.54%  " }
    @{ Row = 40; D = "50.20"; E = "  -1.46%  " }
    @{ Row = 41; D = "434.36"; E = "  +1.53%  " }
    @{ Row = 42; D = "8.75"; E = "  -0.76%  " }
    @{ Row = 43; D = $null; E = "  -0.04%  " }
    @{ Row = 44; D = "2.858.79"; E = "  -1.69%  " }
    @{ Row = 45; D = $null; E = "  -2.02%  " }
    @{ Row = 46; D = $null; E = "  -3.95%  " }
    @{ Row = 47; D = "35.89"; E = "  +2.74%  " }
    @{ Row = 49; D = "124.01"; E = "  +0.52%  " }
    @{ Row = 50; D = $null; E = "  -1.23%  " }
    @{ Row = 51; D = "24.05"; E = "  -2.53%  " }
)

foreach ($u in $updates) {
    if ($null -ne $u.D) {
        $dCell = $ws.Cells.Item($u.Row, 4)
        # Force text storage so numeric-looking prices ("582.48") stay strings
        # (matching the source workbook's inlineStr cells) instead of being
        # auto-converted to numbers by Excel's normal type inference.
        $dCell.NumberFormat = "@"
        $dCell.Value = $u.D
        # Restore the default cell style so no stray formatting is introduced.
        $dCell.Style = "Normal"
    }
    if ($null -ne $u.E) {
        $ws.Cells.Item($u.Row, 5).Value = $u.E
    }
}
